# feat(forms): change default labels
# Update the header row (row 1) of the "verbete" import template with the
# new default field labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "termo"
$ws.Range("E1").Value = "termo superior"
$ws.Range("F1").Value = "formas variantes"
$ws.Range("G1").Value = "termos equivalentes em outros idiomas"
$ws.Range("H1").Value = "ver tambem"
$ws.Range("I1").Value = "fontes"

# Leave the selection where the editor last left it.
$ws.Range("D18").Select()
